$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the current row 956, shifting the existing
# rows 956:978 down to 965:987.
$ws.Rows("956:964").Insert()

# Column layout (row 1 headers):
# A Mercado ID | B Mercado | C Region | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoria ID | J Categoria | K Variedad
# L Calidad | M Volumen | N Precio minimo | O Precio maximo
# P Precio promedio ponderado | Q Unidad de comercializacion | R Origen
# S Precio $/Kg | T Kg / unidad
#
# Columns A, B, C, E, F, G, H, I, J are constant for this entire sheet
# (single market / product / category), so fill them for every new row.

$constA = 6
$constB = "Mercado Mayorista Lo Valledor de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = "Fruta"
$constG = 100103
$constH = "Frutos de hueso (carozo)"
$constI = 100103004
$constJ = "Durazno"

$newRows = @(
    @{ Row = 956; D = 44595; K = "Andross";      L = "Especial"; M = 120; N = 20000;  O = 20000;  P = 20000;  Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 1111; T = 18 },
    @{ Row = 957; D = 44595; K = "Andross";      L = "Primera";  M = 200; N = 16000;  O = 16000;  P = 16000;  Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 889;  T = 18 },
    @{ Row = 958; D = 44595; K = "Andross";      L = "Segunda";  M = 180; N = 12000;  O = 12000;  P = 12000;  Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 667;  T = 18 },
    @{ Row = 959; D = 44595; K = "Carson";       L = "Primera";  M = 8;   N = 280000; O = 280000; P = 280000; Q = "`$/bins (400 kilos)";         R = "Región de O'Higgins"; S = 700;  T = 400 },
    @{ Row = 960; D = 44595; K = "Carson";       L = "Segunda";  M = 14;  N = 250000; O = 250000; P = 250000; Q = "`$/bins (400 kilos)";         R = "Región de O'Higgins"; S = 625;  T = 400 },
    @{ Row = 961; D = 44595; K = "Doctor Davis"; L = "Primera";  M = 10;  N = 250000; O = 250000; P = 250000; Q = "`$/bins (400 kilos)";         R = "Región de O'Higgins"; S = 625;  T = 400 },
    @{ Row = 962; D = 44595; K = "Doctor Davis"; L = "Segunda";  M = 135; N = 14000;  O = 14000;  P = 14000;  Q = "`$/caja 15 kilos granel";     R = "Región de O'Higgins"; S = 933;  T = 15 },
    @{ Row = 963; D = 44595; K = "Doctor Davis"; L = "Segunda";  M = 12;  N = 230000; O = 230000; P = 230000; Q = "`$/bins (400 kilos)";         R = "Región de O'Higgins"; S = 575;  T = 400 },
    @{ Row = 964; D = 44595; K = "Doctor Davis"; L = "Segunda";  M = 135; N = 11000;  O = 11000;  P = 11000;  Q = "`$/caja 15 kilos granel";     R = "Región de O'Higgins"; S = 733;  T = 15 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value2 = $constA
    $ws.Cells.Item($r, 2).Value2 = $constB
    $ws.Cells.Item($r, 3).Value2 = $constC
    $ws.Cells.Item($r, 4).Value2 = $nr.D
    $ws.Cells.Item($r, 5).Value2 = $constE
    $ws.Cells.Item($r, 6).Value2 = $constF
    $ws.Cells.Item($r, 7).Value2 = $constG
    $ws.Cells.Item($r, 8).Value2 = $constH
    $ws.Cells.Item($r, 9).Value2 = $constI
    $ws.Cells.Item($r, 10).Value2 = $constJ
    $ws.Cells.Item($r, 11).Value2 = $nr.K
    $ws.Cells.Item($r, 12).Value2 = $nr.L
    $ws.Cells.Item($r, 13).Value2 = $nr.M
    $ws.Cells.Item($r, 14).Value2 = $nr.N
    $ws.Cells.Item($r, 15).Value2 = $nr.O
    $ws.Cells.Item($r, 16).Value2 = $nr.P
    $ws.Cells.Item($r, 17).Value2 = $nr.Q
    $ws.Cells.Item($r, 18).Value2 = $nr.R
    $ws.Cells.Item($r, 19).Value2 = $nr.S
    $ws.Cells.Item($r, 20).Value2 = $nr.T
}
